$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 2 (pushes the existing rows 2-7 down to 3-8)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "VT_NI_OTH_V4.xlsx"
$ws.Range("B2").Value = "RES"
$ws.Range("C2").Value = "Base"
$ws.Range("D2").Value = "BASE"
$ws.Range("E2").Value = "ncap_life"
$ws.Range("F2").Value = "warning"
$ws.Range("G2").Value = "column is not recognized or invalid"
$ws.Range("H2").Value = "2,8"
$ws.Range("I2").Value = "FI_T"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "-"

# Append a new record row at the end (row 9)
$ws.Range("A9").Value = "VT_SI_OTH_V4.xlsx"
$ws.Range("B9").Value = "RES"
$ws.Range("C9").Value = "Base"
$ws.Range("D9").Value = "BASE"
$ws.Range("E9").Value = "ncap_life"
$ws.Range("F9").Value = "warning"
$ws.Range("G9").Value = "column is not recognized or invalid"
$ws.Range("H9").Value = "2,8"
$ws.Range("I9").Value = "FI_T"
$ws.Range("J9").Value = "-"
$ws.Range("K9").Value = "-"
